$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Diebold-Mariano statistics (DM_Stat, P_Value) and Significativo column
# Row 2 (N_Calib_1=20, N_Calib_2=40)
$ws.Range("C2").Value = 0.3770488950870826
$ws.Range("D2").Value = 0.7084804737371915

# Row 3 (N_Calib_1=20, N_Calib_2=60)
$ws.Range("C3").Value = 0.3661758554726863
$ws.Range("D3").Value = 0.7165026041637743

# Row 4 (N_Calib_1=20, N_Calib_2=100)
$ws.Range("C4").Value = -1.320985395030964
$ws.Range("D4").Value = 0.1953283484710904

# Row 5 (N_Calib_1=20, N_Calib_2=200)
$ws.Range("C5").Value = -1.884892892690875
$ws.Range("D5").Value = 0.06801531879242573
$ws.Range("G5").Value = "No"

# Row 6 (N_Calib_1=40, N_Calib_2=60)
$ws.Range("C6").Value = 0.01839022096354475
$ws.Range("D6").Value = 0.9854350523506523

# Row 7 (N_Calib_1=40, N_Calib_2=100)
$ws.Range("C7").Value = -1.375003981650049
$ws.Range("D7").Value = 0.1781240669822768

# Row 8 (N_Calib_1=40, N_Calib_2=200)
$ws.Range("C8").Value = -2.021922596284605
$ws.Range("D8").Value = 0.05111019185074661
$ws.Range("G8").Value = "No"

# Row 9 (N_Calib_1=60, N_Calib_2=100)
$ws.Range("C9").Value = -1.288110274518346
$ws.Range("D9").Value = 0.2064053868009432

# Row 10 (N_Calib_1=60, N_Calib_2=200)
$ws.Range("C10").Value = -1.89499060913057
$ws.Range("D10").Value = 0.06662294914767442
$ws.Range("G10").Value = "No"

# Row 11 (N_Calib_1=100, N_Calib_2=200)
$ws.Range("C11").Value = -1.586166888251203
$ws.Range("D11").Value = 0.1219594697407995
$ws.Range("G11").Value = "No"
